# correct ID column typo: the "ID" column (G) held plain numbers (1-6 or
# 1-7, repeated per group) but should hold zero-padded text ids like
# "id01".."id07". Re-write each cell's value and give the column a thin
# right-hand border to set it off visually.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idRange = $ws.Range("G2:G14")

# Snapshot the existing numeric values first -- writing the new text
# values into the same cells as we go would otherwise clobber values we
# still need to read on later iterations.
$numbers = @{}
foreach ($cell in $idRange.Cells) {
    $numbers[$cell.Row] = [int]$cell.Value2
}

foreach ($cell in $idRange.Cells) {
    $n = $numbers[$cell.Row]
    $cell.Value = "id{0:D2}" -f $n
}

# Add a thin border on the right edge of the ID column.
$rightBorder = $idRange.Borders.Item(10)  # xlEdgeRight
$rightBorder.LineStyle = 1                # xlContinuous
$rightBorder.Weight = 2                   # xlThin

$ws.Range("E17").Select()
